$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Non-numeric-looking text cells (Coin name, Link, Volume%) ---
$ws.Range("E2").Value = "  +1.05%  "
$ws.Range("E3").Value = "  +0.47%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("E5").Value = "  +3.34%  "
$ws.Range("E6").Value = "  +0.36%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  +2.27%  "
$ws.Range("E9").Value = "  +1.74%  "
$ws.Range("E10").Value = "  +1.98%  "
$ws.Range("E11").Value = "  +0.35%  "
$ws.Range("E12").Value = "  +1.53%  "
$ws.Range("E13").Value = "  +0.75%  "
$ws.Range("E14").Value = "  +1.03%  "
$ws.Range("E15").Value = "  +1.09%  "
$ws.Range("E16").Value = "  +0.96%  "
$ws.Range("E17").Value = "  +1.15%  "
$ws.Range("E18").Value = "  -1.48%  "
$ws.Range("E19").Value = "  +2.11%  "
$ws.Range("E20").Value = "  +1.57%  "
$ws.Range("E21").Value = "  +1.37%  "
$ws.Range("E22").Value = "  +0.03%  "
$ws.Range("E23").Value = "  +1.49%  "
$ws.Range("E24").Value = "  +1.31%  "
$ws.Range("B25").Value = "WrappedeETH"
$ws.Range("C25").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("E25").Value = "  +0.67%  "
$ws.Range("B26").Value = "Kaspa"
$ws.Range("C26").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("E26").Value = "  -0.15%  "
$ws.Range("E27").Value = "  +0.20%  "
$ws.Range("E28").Value = "  +0.54%  "
$ws.Range("E29").Value = "  +0.78%  "
$ws.Range("E30").Value = "  +0.00%  "
$ws.Range("E31").Value = "  +10.53%  "
$ws.Range("E32").Value = "  +0.19%  "
$ws.Range("E33").Value = "  +3.36%  "
$ws.Range("E34").Value = "  -2.80%  "
$ws.Range("E35").Value = "  +5.62%  "
$ws.Range("E36").Value = "  +9.42%  "
$ws.Range("E37").Value = "  +1.06%  "
$ws.Range("E38").Value = "  +2.04%  "
$ws.Range("B39").Value = "Filecoin"
$ws.Range("C39").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("E39").Value = "  +0.91%  "
$ws.Range("B40").Value = "OKB"
$ws.Range("C40").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("E40").Value = "  +1.71%  "
$ws.Range("B41").Value = "Fetch.AI"
$ws.Range("C41").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("E41").Value = "  -0.01%  "
$ws.Range("E42").Value = "  -2.71%  "
$ws.Range("E43").Value = "  +1.42%  "
$ws.Range("E44").Value = "  +1.68%  "
$ws.Range("E45").Value = "  +0.27%  "
$ws.Range("E46").Value = "  +0.07%  "
$ws.Range("E47").Value = "  -0.63%  "
$ws.Range("E48").Value = "  -0.29%  "
$ws.Range("E49").Value = "  +1.55%  "
$ws.Range("E50").Value = "  +0.16%  "
$ws.Range("E51").Value = "  +7.40%  "

# --- Price column cells: force text so numeric-looking strings are not
#     auto-converted to numbers (which would strip formatting like trailing zeros) ---
$priceCells = @("D2", "D3", "D5", "D6", "D7", "D8", "D10", "D11", "D13", "D14", "D15", "D17", "D18", "D19", "D20", "D21", "D23", "D24", "D25", "D26", "D28", "D29", "D31", "D32", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D44", "D45", "D48", "D51")
foreach ($c in $priceCells) { $ws.Range($c).NumberFormat = "@" }
$ws.Range("D2").Value = "60.813.91"
$ws.Range("D3").Value = "2.590.86"
$ws.Range("D5").Value = "522.44"
$ws.Range("D6").Value = "153.89"
$ws.Range("D7").Value = "0.999"
$ws.Range("D8").Value = "0.594"
$ws.Range("D10").Value = "0.105"
$ws.Range("D11").Value = "0.347"
$ws.Range("D13").Value = "3.050.48"
$ws.Range("D14").Value = "60.826.71"
$ws.Range("D15").Value = "21.65"
$ws.Range("D17").Value = "2.599.28"
$ws.Range("D18").Value = "4.74"
$ws.Range("D19").Value = "352.39"
$ws.Range("D20").Value = "10.55"
$ws.Range("D21").Value = "6.21"
$ws.Range("D23").Value = "60.95"
$ws.Range("D24").Value = "0.426"
$ws.Range("D25").Value = "2.716.29"
$ws.Range("D26").Value = "0.166"
$ws.Range("D28").Value = "0.0₃0847"
$ws.Range("D29").Value = "7.35"
$ws.Range("D31").Value = "6.33"
$ws.Range("D32").Value = "19.32"
$ws.Range("D34").Value = "148.84"
$ws.Range("D35").Value = "4.20"
$ws.Range("D36").Value = "0.940"
$ws.Range("D37").Value = "1.20"
$ws.Range("D38").Value = "1.49"
$ws.Range("D39").Value = "3.79"
$ws.Range("D40").Value = "36.45"
$ws.Range("D41").Value = "0.847"
$ws.Range("D42").Value = "286.91"
$ws.Range("D44").Value = "0.623"
$ws.Range("D45").Value = "0.0560"
$ws.Range("D48").Value = "4.87"
$ws.Range("D51").Value = "18.99"
foreach ($c in $priceCells) { $ws.Range($c).Style = "Normal" }
